# Textgrid Speak pts 40-41, Praat output files
#
# Adds two new participant rows (41 and 42, i.e. spreadsheet rows 43 and 44)
# to the "Participants" sheet: Angelika Botero Montaña and Alejandro Muñoz
# Pérez, including their formatting (copied from existing fully-populated
# rows), notes, dates of testing/birth and sex, plus a couple of cosmetic
# view/page-setup tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participants")

# --- Formatting ----------------------------------------------------------
# Row 43 needs the same per-column styles as row 30 (A,B,C,D,E,F,G,H,I,J,K
# already line up exactly with what row 43 should become).
$ws.Range("A30:K30").Copy() | Out-Null
$ws.Range("A43:K43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 44 matches row 39's per-column styles for every column except I
# (Notes), so copy row 39 first and then patch column I separately.
$ws.Range("A39:K39").Copy() | Out-Null
$ws.Range("A44:K44").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I30").Copy() | Out-Null
$ws.Range("I44").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data ------------------------------------------------------------------
# Row 43 - participant 42: Angelika Botero Montaña
$ws.Range("B43").Value = "Angelika Botero Montaña"
$ws.Range("C43").Value = 10275
$ws.Range("H43").Value = 45307
$ws.Range("I43").Value = "Used words from South American Spanish in the familliarisation."
$ws.Range("J43").Value = 33061
$ws.Range("K43").Value = "NR"

# Row 44 - participant 43: Alejandro Muñoz Pérez
$ws.Range("B44").Value = "Alejandro Muñoz Pérez"
$ws.Range("C44").Value = 6835
$ws.Range("H44").Value = 45309
$ws.Range("I44").Value = "In Practice2 I wrote his code ID as 430201 instead of 430102."
$ws.Range("J44").Value = 33718
$ws.Range("K44").Value = "M"

# --- Page setup --------------------------------------------------------
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait

# --- View state (best effort) -------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C36").Select()
